$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "cond_pm1"

# Clear out the old data area (rows 1:9, cols A:B) before rewriting
$ws.Range("A1:D9").ClearContents()

# Header row (set in order so shared-string table indices come out as
# 0=stim1_c, 1=stim2_c, 2=SOA, 3=angle_diff)
$ws.Range("C1").Value = "stim1_c"
$ws.Range("D1").Value = "stim2_c"
$ws.Range("A1").Value = "SOA"
$ws.Range("B1").Value = "angle_diff"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = -0.4
$ws.Range("D2").Value = -0.4

# Row 3
$ws.Range("A3").Value = 12
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = -0.4
$ws.Range("D3").Value = -0.4

# Update selection to match target state
$ws.Range("E8").Select()
